$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.463.69'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '1.908.85'
$ws.Range("E3").Value = '  +0.01%  '
$ws.Range("E4").Value = '  +0.69%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.78'
$ws.Range("E5").Value = '  +0.53%  '
$ws.Range("E6").Value = '  +0.53%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4838'
$ws.Range("E7").Value = '  +2.30%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4061'
$ws.Range("E8").Value = '  +0.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08165'
$ws.Range("E9").Value = '  +1.78%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.012'
$ws.Range("E10").Value = '  +1.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.50'
$ws.Range("E11").Value = '  +3.54%  '
$ws.Range("D12").Value = '1.908.33'
$ws.Range("E12").Value = '  -1.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.034'
$ws.Range("E13").Value = '  +2.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.185'
$ws.Range("E14").Value = '  +1.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.46'
$ws.Range("E15").Value = '  +1.19%  '
$ws.Range("E16").Value = '  +0.69%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.06761'
$ws.Range("E17").Value = '  +1.88%  '
$ws.Range("E18").Value = '  +0.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.69'
$ws.Range("E19").Value = '  +0.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.006'
$ws.Range("E20").Value = '  +0.51%  '
$ws.Range("D21").Value = '29.499.80'
$ws.Range("E21").Value = '  +0.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.628'
$ws.Range("E22").Value = '  +1.99%  '
$ws.Range("E23").Value = '  +2.47%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.190'
$ws.Range("E24").Value = '  -0.51%  '
$ws.Range("D25").Value = '2.160.86'
$ws.Range("E25").Value = '  +0.28%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.75'
$ws.Range("E26").Value = '  +1.91%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.503'
$ws.Range("E27").Value = '  +7.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.08'
$ws.Range("E28").Value = '  +1.59%  '
$ws.Range("E29").Value = '  +0.97%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.41'
$ws.Range("E30").Value = '  +2.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.028'
$ws.Range("E31").Value = '  -3.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09525'
$ws.Range("E32").Value = '  +0.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.517'
$ws.Range("E33").Value = '  +2.76%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.560'
$ws.Range("E34").Value = '  +0.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.390'
$ws.Range("E35").Value = '  -1.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02276'
$ws.Range("E36").Value = '  +1.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06118'
$ws.Range("E37").Value = '  +0.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.179'
$ws.Range("E38").Value = '  +0.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.83'
$ws.Range("E39").Value = '  +7.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5958'
$ws.Range("E40").Value = '  +2.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.971'
$ws.Range("E41").Value = '  -1.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1857'
$ws.Range("E42").Value = '  +1.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.277'
$ws.Range("E43").Value = '  +0.43%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.377'
$ws.Range("E44").Value = '  -5.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.48'
$ws.Range("E45").Value = '  +2.97%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07601'
$ws.Range("E46").Value = '  -2.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5573'
$ws.Range("E47").Value = '  +1.42%  '
$ws.Range("E48").Value = '  +2.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '116.69'
$ws.Range("E49").Value = '  +2.79%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.58'
$ws.Range("E50").Value = '  +1.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.410'
$ws.Range("E51").Value = '  +2.85%  '
